$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 30.87085333333333
$ws.Range("H2").Value = 92.61256
$ws.Range("I2").Value = 0.2985789950947061
$ws.Range("J2").Value = 0.2985789950947061
$ws.Range("M2").Value = 0.1994653333333334
$ws.Range("N2").Value = 0.598396
$ws.Range("O2").Value = 0.01676579960230272
$ws.Range("P2").Value = 0.01676579960230271
$ws.Range("Q2").Value = 6.157665050417778
$ws.Range("R2").Value = 55.41898545376
$ws.Range("S2").Value = 0.005005915597214768
$ws.Range("T2").Value = 0.005005915597214768
$ws.Range("G3").Value = 30.87085333333333
$ws.Range("H3").Value = 92.61256
$ws.Range("I3").Value = 0.2985789950947061
$ws.Range("J3").Value = 0.2985789950947061
$ws.Range("O3").Value = 0.03203779682023726
$ws.Range("P3").Value = 0.03203779682023726
$ws.Range("Q3").Value = 11.76669329539556
$ws.Range("R3").Value = 105.90023965856
$ws.Range("S3").Value = 0.009565813179634812
$ws.Range("T3").Value = 0.009565813179634813
$ws.Range("G4").Value = 30.87085333333333
$ws.Range("H4").Value = 92.61256
$ws.Range("I4").Value = 0.2985789950947061
$ws.Range("J4").Value = 0.2985789950947061
$ws.Range("M4").Value = 0.2888043333333333
$ws.Range("N4").Value = 0.8664129999999999
$ws.Range("O4").Value = 0.02427507324719734
$ws.Range("P4").Value = 0.02427507324719734
$ws.Range("Q4").Value = 8.915636216364442
$ws.Range("R4").Value = 80.24072594727998
$ws.Range("S4").Value = 0.007248026975998565
$ws.Range("T4").Value = 0.007248026975998566
$ws.Range("G5").Value = 30.87085333333333
$ws.Range("H5").Value = 92.61256
$ws.Range("I5").Value = 0.2985789950947061
$ws.Range("J5").Value = 0.2985789950947061
$ws.Range("M5").Value = 11.02772766666667
$ws.Range("N5").Value = 33.083183
$ws.Range("O5").Value = 0.9269213303302627
$ws.Range("P5").Value = 0.9269213303302626
$ws.Range("Q5").Value = 340.4353633976089
$ws.Range("R5").Value = 3063.91827057848
$ws.Range("S5").Value = 0.2767592393418579
$ws.Range("T5").Value = 0.276759239341858
$ws.Range("G6").Value = 33.793597
$ws.Range("I6").Value = 0.3268474027571036
$ws.Range("J6").Value = 0.3268474027571037
$ws.Range("M6").Value = 0.1994653333333334
$ws.Range("N6").Value = 0.598396
$ws.Range("O6").Value = 0.01676579960230272
$ws.Range("P6").Value = 0.01676579960230271
$ws.Range("Q6").Value = 6.740651090137334
$ws.Range("R6").Value = 60.665859811236
$ws.Range("S6").Value = 0.005479858055158724
$ws.Range("T6").Value = 0.005479858055158723
$ws.Range("G7").Value = 33.793597
$ws.Range("I7").Value = 0.3268474027571036
$ws.Range("J7").Value = 0.3268474027571037
$ws.Range("O7").Value = 0.03203779682023726
$ws.Range("P7").Value = 0.03203779682023726
$ws.Range("Q7").Value = 12.88072237439066
$ws.Range("S7").Value = 0.01047147068075434
$ws.Range("T7").Value = 0.01047147068075434
$ws.Range("G8").Value = 33.793597
$ws.Range("I8").Value = 0.3268474027571036
$ws.Range("J8").Value = 0.3268474027571037
$ws.Range("M8").Value = 0.2888043333333333
$ws.Range("N8").Value = 0.8664129999999999
$ws.Range("O8").Value = 0.02427507324719734
$ws.Range("P8").Value = 0.02427507324719734
$ws.Range("Q8").Value = 9.759737252520331
$ws.Range("R8").Value = 87.83763527268299
$ws.Range("S8").Value = 0.0079342446425849
$ws.Range("T8").Value = 0.0079342446425849
$ws.Range("G9").Value = 33.793597
$ws.Range("I9").Value = 0.3268474027571036
$ws.Range("J9").Value = 0.3268474027571037
$ws.Range("M9").Value = 11.02772766666667
$ws.Range("N9").Value = 33.083183
$ws.Range("O9").Value = 0.9269213303302627
$ws.Range("P9").Value = 0.9269213303302626
$ws.Range("Q9").Value = 372.6665845930836
$ws.Range("R9").Value = 3353.999261337753
$ws.Range("S9").Value = 0.3029618293786057
$ws.Range("T9").Value = 0.3029618293786057
$ws.Range("G10").Value = 2.981185666666667
$ws.Range("H10").Value = 8.943557
$ws.Range("I10").Value = 0.02883365130639111
$ws.Range("J10").Value = 0.02883365130639111
$ws.Range("M10").Value = 0.1994653333333334
$ws.Range("N10").Value = 0.598396
$ws.Range("O10").Value = 0.01676579960230272
$ws.Range("P10").Value = 0.01676579960230271
$ws.Range("Q10").Value = 0.5946431927302223
$ws.Range("R10").Value = 5.351788734572001
$ws.Range("S10").Value = 0.0004834192196056272
$ws.Range("T10").Value = 0.0004834192196056271
$ws.Range("G11").Value = 2.981185666666667
$ws.Range("H11").Value = 8.943557
$ws.Range("I11").Value = 0.02883365130639111
$ws.Range("J11").Value = 0.02883365130639111
$ws.Range("O11").Value = 0.03203779682023726
$ws.Range("P11").Value = 0.03203779682023726
$ws.Range("Q11").Value = 1.136304753792444
$ws.Range("R11").Value = 10.226742784132
$ws.Range("S11").Value = 0.0009237666621397269
$ws.Range("T11").Value = 0.0009237666621397269
$ws.Range("G12").Value = 2.981185666666667
$ws.Range("H12").Value = 8.943557
$ws.Range("I12").Value = 0.02883365130639111
$ws.Range("J12").Value = 0.02883365130639111
$ws.Range("M12").Value = 0.2888043333333333
$ws.Range("N12").Value = 0.8664129999999999
$ws.Range("O12").Value = 0.02427507324719734
$ws.Range("P12").Value = 0.02427507324719734
$ws.Range("Q12").Value = 0.8609793390045554
$ws.Range("R12").Value = 7.748814051040999
$ws.Range("S12").Value = 0.0006999389974467914
$ws.Range("T12").Value = 0.0006999389974467913
$ws.Range("G13").Value = 2.981185666666667
$ws.Range("H13").Value = 8.943557
$ws.Range("I13").Value = 0.02883365130639111
$ws.Range("J13").Value = 0.02883365130639111
$ws.Range("M13").Value = 11.02772766666667
$ws.Range("N13").Value = 33.083183
$ws.Range("O13").Value = 0.9269213303302627
$ws.Range("P13").Value = 0.9269213303302626
$ws.Range("Q13").Value = 32.87570365577011
$ws.Range("R13").Value = 295.881332901931
$ws.Range("S13").Value = 0.02672652642719896
$ws.Range("T13").Value = 0.02672652642719896
$ws.Range("G14").Value = 35.74694633333333
$ws.Range("H14").Value = 107.240839
$ws.Range("I14").Value = 0.3457399508417991
$ws.Range("J14").Value = 0.3457399508417991
$ws.Range("M14").Value = 0.1994653333333334
$ws.Range("N14").Value = 0.598396
$ws.Range("O14").Value = 0.01676579960230272
$ws.Range("P14").Value = 0.01676579960230271
$ws.Range("Q14").Value = 7.130276566027112
$ws.Range("R14").Value = 64.17248909424401
$ws.Range("S14").Value = 0.005796606730323596
$ws.Range("T14").Value = 0.005796606730323594
$ws.Range("G15").Value = 35.74694633333333
$ws.Range("H15").Value = 107.240839
$ws.Range("I15").Value = 0.3457399508417991
$ws.Range("J15").Value = 0.3457399508417991
$ws.Range("O15").Value = 0.03203779682023726
$ws.Range("P15").Value = 0.03203779682023726
$ws.Range("Q15").Value = 13.62525840181822
$ws.Range("R15").Value = 122.627325616364
$ws.Range("S15").Value = 0.01107674629770838
$ws.Range("T15").Value = 0.01107674629770838
$ws.Range("G16").Value = 35.74694633333333
$ws.Range("H16").Value = 107.240839
$ws.Range("I16").Value = 0.3457399508417991
$ws.Range("J16").Value = 0.3457399508417991
$ws.Range("M16").Value = 0.2888043333333333
$ws.Range("N16").Value = 0.8664129999999999
$ws.Range("O16").Value = 0.02427507324719734
$ws.Range("P16").Value = 0.02427507324719734
$ws.Range("Q16").Value = 10.32387300450078
$ws.Range("R16").Value = 92.91485704050699
$ws.Range("S16").Value = 0.008392862631167081
$ws.Range("T16").Value = 0.00839286263116708
$ws.Range("G17").Value = 35.74694633333333
$ws.Range("H17").Value = 107.240839
$ws.Range("I17").Value = 0.3457399508417991
$ws.Range("J17").Value = 0.3457399508417991
$ws.Range("M17").Value = 11.02772766666667
$ws.Range("N17").Value = 33.083183
$ws.Range("O17").Value = 0.9269213303302627
$ws.Range("P17").Value = 0.9269213303302626
$ws.Range("Q17").Value = 394.2075890789485
$ws.Range("R17").Value = 3547.868301710537
$ws.Range("S17").Value = 0.3204737351826
$ws.Range("T17").Value = 0.3204737351826
